$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_Beton@Erde"
$ws.Range("B2").Value = "Beton@Erde"
$ws.Range("C2").Value = 45049
$ws.Range("D2").Value = 45050

# Update row 3
$ws.Range("A3").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_Beton@Schalen"
$ws.Range("B3").Value = "Beton@Schalen"
$ws.Range("C3").Value = 45050
$ws.Range("D3").Value = 45051

# Add new row 4
$ws.Range("A4").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_Beton@Vorbereitung"
$ws.Range("B4").Value = "Beton@Vorbereitung"
$ws.Range("C4").Value = 45048
$ws.Range("D4").Value = 45049

# Apply the same number format as C3/D3 (date format) to C4/D4
$ws.Range("C4:D4").NumberFormat = $ws.Range("C3:D3").NumberFormat
